# Auto-generated Excel COM-interop script applying the scheduled-runner value updates
# described by the upstream diff for Halicarnassus_Profits.xlsx.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 8238.299999999999
$ws.Range("I40").Value = 5971.75
$ws.Range("J40").Value = 9749.333000000001
$ws.Range("K40").Value = 5971.75
$ws.Range("L40").Value = 9749.333000000001
$ws.Range("M40").Value = -5796.75
$ws.Range("N40").Value = -10099.333
$ws.Range("H43").Value = 2724.5
$ws.Range("I43").Value = 3000
$ws.Range("J43").Value = 2449
$ws.Range("K43").Value = 3000
$ws.Range("L43").Value = 2449
$ws.Range("M43").Value = -2931
$ws.Range("N43").Value = -2587
$ws.Range("H53").Value = 388.8125
$ws.Range("I53").Value = 516.7778
$ws.Range("J53").Value = 224.28572
$ws.Range("K53").Value = 516.7778
$ws.Range("L53").Value = 224.28572
$ws.Range("M53").Value = 120.2222
$ws.Range("N53").Value = -1498.28572
$ws.Range("H76").Value = 1879.8
$ws.Range("I76").Value = 1879.8
$ws.Range("K76").Value = 1879.8
$ws.Range("M76").Value = -1564.8
$ws.Range("H79").Value = 1879.8
$ws.Range("I79").Value = 1879.8
$ws.Range("K79").Value = 1879.8
$ws.Range("M79").Value = -787.8
$ws.Range("H86").Value = 1000
$ws.Range("I86").Value = 1000
$ws.Range("K86").Value = 1000
$ws.Range("M86").Value = 123
$ws.Range("H89").Value = 1000
$ws.Range("I89").Value = 1000
$ws.Range("K89").Value = 5000
$ws.Range("M89").Value = 616
$ws.Range("H106").Value = 1000
$ws.Range("I106").Value = 1000
$ws.Range("K106").Value = 1000
$ws.Range("M106").Value = -369
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H116").Value = 4717.4
$ws.Range("J116").Value = 4694.5
$ws.Range("L116").Value = 4694.5
$ws.Range("N116").Value = -11578.5
$ws.Range("H129").Value = 1998.3334
$ws.Range("I129").Value = 1499.5
$ws.Range("K129").Value = 4498.5
$ws.Range("M129").Value = 501.5
$ws.Range("H137").Value = 4299.8945
$ws.Range("I137").Value = 919.1539
$ws.Range("J137").Value = 11624.833
$ws.Range("K137").Value = 2757.4617
$ws.Range("L137").Value = 34874.499
$ws.Range("M137").Value = -207.4616999999998
$ws.Range("N137").Value = -39974.499
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1296.4286
$ws.Range("I2").Value = 1296.4286
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1296.4286
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -1183.4286
$ws.Range("N2").ClearContents()
$ws.Range("H45").Value = 3292.4614
$ws.Range("I45").Value = 2112.8333
$ws.Range("K45").Value = 2112.8333
$ws.Range("M45").Value = -1735.8333
$ws.Range("H63").Value = 5458.4287
$ws.Range("I63").Value = 2851
$ws.Range("K63").Value = 2851
$ws.Range("M63").Value = -2165
$ws.Range("H66").Value = 5458.4287
$ws.Range("I66").Value = 2851
$ws.Range("K66").Value = 14255
$ws.Range("M66").Value = -10823
$ws.Range("H74").Value = 3337.6155
$ws.Range("I74").Value = 2532.1
$ws.Range("K74").Value = 2532.1
$ws.Range("M74").Value = -1658.1
$ws.Range("H77").Value = 3337.6155
$ws.Range("I77").Value = 2532.1
$ws.Range("K77").Value = 12660.5
$ws.Range("M77").Value = -8292.5
$ws.Range("H116").Value = 1296.4286
$ws.Range("I116").Value = 1296.4286
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 1296.4286
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 997.5714
$ws.Range("N116").ClearContents()
$ws.Range("H117").Value = 25000
$ws.Range("J117").Value = 25000
$ws.Range("L117").Value = 25000
$ws.Range("N117").Value = -34178
$ws.Range("H119").Value = 20000
$ws.Range("J119").Value = 20000
$ws.Range("L119").Value = 20000
$ws.Range("N119").Value = -29676
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1296.4286
$ws.Range("I3").Value = 1296.4286
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1296.4286
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -1182.4286
$ws.Range("N3").ClearContents()
$ws.Range("H94").Value = 909.0714
$ws.Range("I94").Value = 677.3333
$ws.Range("K94").Value = 677.3333
$ws.Range("M94").Value = -226.3333
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 115
$ws.Range("I7").Value = 67.333336
$ws.Range("J7").Value = 150.75
$ws.Range("K7").Value = 67.333336
$ws.Range("L7").Value = 150.75
$ws.Range("M7").Value = 45.666664
$ws.Range("N7").Value = -376.75
$ws.Range("H16").Value = 1380
$ws.Range("I16").Value = 10
$ws.Range("J16").Value = 2750
$ws.Range("K16").Value = 10
$ws.Range("L16").Value = 2750
$ws.Range("M16").Value = 277
$ws.Range("N16").Value = -3324
$ws.Range("H113").Value = 1380
$ws.Range("I113").Value = 10
$ws.Range("J113").Value = 2750
$ws.Range("K113").Value = 10
$ws.Range("L113").Value = 2750
$ws.Range("M113").Value = 2160
$ws.Range("N113").Value = -7090
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1123.7307
$ws.Range("I5").Value = 1256.5454
$ws.Range("J5").Value = 1026.3334
$ws.Range("K5").Value = 3769.6362
$ws.Range("L5").Value = 3079.0002
$ws.Range("M5").Value = -3657.6362
$ws.Range("N5").Value = -3303.0002
$ws.Range("H92").Value = 691.75
$ws.Range("I92").Value = 580.1
$ws.Range("K92").Value = 1740.3
$ws.Range("M92").Value = -492.3000000000002
$ws.Range("H115").Value = 2499.5
$ws.Range("I115").Value = 2000
$ws.Range("K115").Value = 6000
$ws.Range("M115").Value = -4825
$ws.Range("H131").Value = 1657.6
$ws.Range("J131").Value = 1896.5
$ws.Range("L131").Value = 5689.5
$ws.Range("N131").Value = -15769.5
$ws.Range("H135").Value = 1123.7307
$ws.Range("I135").Value = 1256.5454
$ws.Range("J135").Value = 1026.3334
$ws.Range("K135").Value = 11308.9086
$ws.Range("L135").Value = 9237.000599999999
$ws.Range("M135").Value = -8773.908599999999
$ws.Range("N135").Value = -14307.0006
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H70").Value = 3448.75
$ws.Range("I70").Value = 3448.75
$ws.Range("K70").Value = 3448.75
$ws.Range("M70").Value = -3178.75
$ws.Range("H73").Value = 3448.75
$ws.Range("I73").Value = 3448.75
$ws.Range("K73").Value = 3448.75
$ws.Range("M73").Value = -2512.75
$ws.Range("H132").Value = 1000012
$ws.Range("I132").Value = 1000012
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3000036
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2997506
$ws.Range("N132").ClearContents()
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 904
$ws.Range("I7").Value = 904
$ws.Range("K7").Value = 904
$ws.Range("M7").Value = -792
$ws.Range("H16").Value = 2333.3333
$ws.Range("I16").Value = 2500
$ws.Range("K16").Value = 2500
$ws.Range("M16").Value = -2330
$ws.Range("H43").Value = 28999.5
$ws.Range("J43").Value = 28999.5
$ws.Range("L43").Value = 28999.5
$ws.Range("N43").Value = -29385.5
$ws.Range("H126").Value = 904
$ws.Range("I126").Value = 904
$ws.Range("K126").Value = 2712
$ws.Range("M126").Value = -242
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1468.5
$ws.Range("I100").Value = 1468.5
$ws.Range("K100").Value = 2937
$ws.Range("M100").Value = -2396
$ws.Range("H107").Value = 794.5454999999999
$ws.Range("I107").Value = 819.4286
$ws.Range("J107").Value = 751
$ws.Range("K107").Value = 2458.2858
$ws.Range("L107").Value = 2253
$ws.Range("M107").Value = -538.2857999999997
$ws.Range("N107").Value = -6093
